# Add new column 'Servised by' to Card24
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# New column goes right after the existing "Correction" column (N), i.e. column O (15)
# Mirror column N's header formatting (bold, bordered, centered/top-aligned) onto O1
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Servised by"

# Extend the new column down through the existing data rows (2-12) so the
# sheet's used range covers the new column for every row, same as the rest
# of the table. These cells stay blank - same as the other (empty) source data.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 15).Value = " "
}
